$d = $word.ActiveDocument

# Both title paragraphs get a yellow (FFFF00) font color applied to the
# run text AND to the paragraph mark itself (so the pilcrow / end-of-
# paragraph run properties carry the same color), plus a bright-green
# highlight on the run text.
$titles = @("Edición y operaciones básicas con el texto", "TEXTO RANDOM")

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.Trim()
    foreach ($title in $titles) {
        if ($text -eq $title) {
            $rng = $para.Range
            $rng.Font.Color = 65535
            $rng.HighlightColorIndex = 4
        }
    }
}
